# Add new columns I (I0) and J (IF) with header + data for rows 2..29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..29 : row -> @(I value, J value)
$data = @{
    2  = @(8, 9)
    3  = @(2, 3)
    4  = @(3, 5)
    5  = @(3, 4)
    6  = @(1, 2)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 1)
    11 = @(1, 4)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(1, 7)
    16 = @(1, 6)
    17 = @(1, 5)
    18 = @(1, 5)
    19 = @(2, 8)
    20 = @(1, 4)
    21 = @(1, 5)
    22 = @(1, 4)
    23 = @(1, 5)
    24 = @(1, 5)
    25 = @(1, 4)
    26 = @(1, 6)
    27 = @(1, 4)
    28 = @(5, 7)
    29 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
